$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update EC database: swap the two account rows (16 and 18) ---
# Row 16 becomes: 45371686 / LUCIA ISABEL ACEVEDO TORRES / 1711 / 14755 / 877803
# Row 18 becomes: 1143348570 / HECTOR GUILLERMO ARRIETA RODELO / 1802 / 15625 / 961420
$ws.Range("C16").Value = "45371686"
$ws.Range("D16").Value = "LUCIA ISABEL ACEVEDO TORRES"
$ws.Range("E16").Value = "1711"
$ws.Range("F16").Value = 14755
$ws.Range("G16").Value = 877803

$ws.Range("C18").Value = "1143348570"
$ws.Range("D18").Value = "HECTOR GUILLERMO ARRIETA RODELO"
$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 15625
$ws.Range("G18").Value = 961420

# --- Reposition the logo image (moved left by 13.5pt / 171450 EMU) ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 59.090551181102356
$shp.Top = 19.405511811023622
$shp.Width = 76.81889763779527
$shp.Height = 48.188976377952756

Write-Output "done"
